$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values to reflect repulled data / recalculated means
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -12
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -15
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -4
$ws.Range("F10").Value = -2
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = -9
$ws.Range("F15").Value = -3
$ws.Range("F16").Value = -6
$ws.Range("F17").Value = 1
